# Apply updated mounted-pipeline values (Step1_Data raw signal, Step2_Sj
# cumulative sums, and the derived Step3_DataPts_* threshold-crossing stats)
# per the authoritative diff.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Step1_Data")
$ws1.Range("AN3").Value = 0
$ws1.Range("AO3").Value = 0.0857002670956465
$ws1.Range("AP3").Value = 0.1890461935933762
$ws1.Range("AQ3").Value = 0.04956014656443426
$ws1.Range("AR3").Value = 0.0001989547730391222
$ws1.Range("AS3").Value = 0.02590736071296167
$ws1.Range("AT3").Value = 0.03595149208127125
$ws1.Range("AU3").Value = 0.02910088577780513
$ws1.Range("AV3").Value = 0.00005132133393345735
$ws1.Range("AW3").Value = 0.1415687674698234
$ws1.Range("AX3").Value = 0.0002792704462042901
$ws1.Range("AY3").Value = 0.01575130443425728
$ws1.Range("AZ3").Value = 0.03145485226612948
$ws1.Range("BA3").Value = 0.05256234912724426
$ws1.Range("BB3").Value = 0.002676939309930768
$ws1.Range("BC3").Value = 0.06769425923342286
$ws1.Range("BD3").Value = 0.04980413067475936
$ws1.Range("BE3").Value = 0.006510181166959709
$ws1.Range("BF3").Value = 0.03548029983793137
$ws1.Range("BG3").Value = 0.02558535817866346
$ws1.Range("BH3").Value = 0.07236395926689111
$ws1.Range("BI3").Value = 0.01023901866812264
$ws1.Range("BJ3").Value = 0.0007662616005865862
$ws1.Range("BK3").Value = 0.004021469802913399
$ws1.Range("BL3").Value = 0.009994359778589416
$ws1.Range("BM3").Value = 0.006917657032937319
$ws1.Range("BN3").Value = 0.0005727799414455385
$ws1.Range("BO3").Value = 0.0003441378246944437
$ws1.Range("BP3").Value = 0.002845595931612141
$ws1.Range("BQ3").Value = 0.01314992227945014
$ws1.Range("BR3").Value = 0.01426028977738931
$ws1.Range("BS3").Value = 0.006634837664965398
$ws1.Range("BT3").Value = 0.0130053763526089
$ws1.Range("AN5").Value = 0
$ws1.Range("AO5").Value = 0.1410612498598294
$ws1.Range("AP5").Value = 0.2181479010517564
$ws1.Range("AQ5").Value = 0.139009907792523
$ws1.Range("AR5").Value = 0.0009527868245702934
$ws1.Range("AS5").Value = 0.02004298887845229
$ws1.Range("AT5").Value = 0.01058938382466851
$ws1.Range("AU5").Value = 0.07119636642925078
$ws1.Range("AV5").Value = 0.0007846984720188629
$ws1.Range("AW5").Value = 0.06062452040204066
$ws1.Range("AX5").Value = 0.001648425688132459
$ws1.Range("AY5").Value = 0.01979438622230643
$ws1.Range("AZ5").Value = 0.0004095671879848328
$ws1.Range("BA5").Value = 0.01138671932311719
$ws1.Range("BB5").Value = 0.001017114295013625
$ws1.Range("BC5").Value = 0.01411424735657834
$ws1.Range("BD5").Value = 0.02783086720667993
$ws1.Range("BE5").Value = 0.003574487511241921
$ws1.Range("BF5").Value = 0.006432554633495169
$ws1.Range("BG5").Value = 0.01967492750313842
$ws1.Range("BH5").Value = 0.07855282177717011
$ws1.Range("BI5").Value = 0.007451617174862648
$ws1.Range("BJ5").Value = 0.006053459778593355
$ws1.Range("BK5").Value = 0.00454204377325507
$ws1.Range("BL5").Value = 0.008157051898241114
$ws1.Range("BM5").Value = 0.01545057873879046
$ws1.Range("BN5").Value = 0.0002816753304646522
$ws1.Range("BO5").Value = 0.0002995345903038083
$ws1.Range("BP5").Value = 0.01399543170076312
$ws1.Range("BQ5").Value = 0.03781271184711886
$ws1.Range("BR5").Value = 0.04060254888510798
$ws1.Range("BS5").Value = 0.007166029367500628
$ws1.Range("BT5").Value = 0.01134139467502989

$ws2 = $wb.Worksheets.Item("Step2_Sj")
$ws2.Range("AN3").Value = 0
$ws2.Range("AO3").Value = 0.0857002670956465
$ws2.Range("AP3").Value = 0.2747464606890226
$ws2.Range("AQ3").Value = 0.3243066072534569
$ws2.Range("AR3").Value = 0.324505562026496
$ws2.Range("AS3").Value = 0.3504129227394577
$ws2.Range("AT3").Value = 0.3863644148207289
$ws2.Range("AU3").Value = 0.415465300598534
$ws2.Range("AV3").Value = 0.4155166219324675
$ws2.Range("AW3").Value = 0.5570853894022909
$ws2.Range("AX3").Value = 0.5573646598484951
$ws2.Range("AY3").Value = 0.5731159642827525
$ws2.Range("AZ3").Value = 0.6045708165488819
$ws2.Range("BA3").Value = 0.6571331656761261
$ws2.Range("BB3").Value = 0.6598101049860569
$ws2.Range("BC3").Value = 0.7275043642194797
$ws2.Range("BD3").Value = 0.777308494894239
$ws2.Range("BE3").Value = 0.7838186760611987
$ws2.Range("BF3").Value = 0.81929897589913
$ws2.Range("BG3").Value = 0.8448843340777935
$ws2.Range("BH3").Value = 0.9172482933446846
$ws2.Range("BI3").Value = 0.9274873120128072
$ws2.Range("BJ3").Value = 0.9282535736133938
$ws2.Range("BK3").Value = 0.9322750434163072
$ws2.Range("BL3").Value = 0.9422694031948966
$ws2.Range("BM3").Value = 0.9491870602278339
$ws2.Range("BN3").Value = 0.9497598401692795
$ws2.Range("BO3").Value = 0.9501039779939739
$ws2.Range("BP3").Value = 0.9529495739255861
$ws2.Range("BQ3").Value = 0.9660994962050362
$ws2.Range("BR3").Value = 0.9803597859824256
$ws2.Range("BS3").Value = 0.986994623647391
$ws2.Range("BT3").Value = 0.9999999999999999
$ws2.Range("AN5").Value = 0
$ws2.Range("AO5").Value = 0.1410612498598294
$ws2.Range("AP5").Value = 0.3592091509115858
$ws2.Range("AQ5").Value = 0.4982190587041088
$ws2.Range("AR5").Value = 0.499171845528679
$ws2.Range("AS5").Value = 0.5192148344071313
$ws2.Range("AT5").Value = 0.5298042182317999
$ws2.Range("AU5").Value = 0.6010005846610507
$ws2.Range("AV5").Value = 0.6017852831330696
$ws2.Range("AW5").Value = 0.6624098035351103
$ws2.Range("AX5").Value = 0.6640582292232428
$ws2.Range("AY5").Value = 0.6838526154455492
$ws2.Range("AZ5").Value = 0.684262182633534
$ws2.Range("BA5").Value = 0.6956489019566512
$ws2.Range("BB5").Value = 0.6966660162516648
$ws2.Range("BC5").Value = 0.7107802636082431
$ws2.Range("BD5").Value = 0.738611130814923
$ws2.Range("BE5").Value = 0.7421856183261649
$ws2.Range("BF5").Value = 0.74861817295966
$ws2.Range("BG5").Value = 0.7682931004627984
$ws2.Range("BH5").Value = 0.8468459222399686
$ws2.Range("BI5").Value = 0.8542975394148312
$ws2.Range("BJ5").Value = 0.8603509991934246
$ws2.Range("BK5").Value = 0.8648930429666797
$ws2.Range("BL5").Value = 0.8730500948649208
$ws2.Range("BM5").Value = 0.8885006736037112
$ws2.Range("BN5").Value = 0.8887823489341758
$ws2.Range("BO5").Value = 0.8890818835244797
$ws2.Range("BP5").Value = 0.9030773152252428
$ws2.Range("BQ5").Value = 0.9408900270723617
$ws2.Range("BR5").Value = 0.9814925759574696
$ws2.Range("BS5").Value = 0.9886586053249703

$ws3 = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws3.Range("F3").Value = 0.5570853894022909
$ws3.Range("D5").Value = 44
$ws3.Range("F5").Value = 0.5192148344071313
$ws3.Range("G5").Value = 6

$ws4 = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws4.Range("F3").Value = 0.7275043642194797
$ws4.Range("D5").Value = 54
$ws4.Range("F5").Value = 0.7107802636082431
$ws4.Range("G5").Value = 16

$ws5 = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws5.Range("F3").Value = 0.81929897589913
$ws5.Range("F5").Value = 0.8468459222399686

$ws6 = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws6.Range("F3").Value = 0.9172482933446846
$ws6.Range("D5").Value = 67
$ws6.Range("F5").Value = 0.9030773152252428
$ws6.Range("G5").Value = 29
